# Housecleaning after abandoning the QFN16 SI5351 package
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Status" column (C) entries that referred to sourcing notes
# that are no longer relevant, except C4 ("Future") which stays.
$ws.Range("C2").Clear()
$ws.Range("C3").Clear()
$ws.Range("C5").Clear()
$ws.Range("C6").Clear()
$ws.Range("C11").Clear()
$ws.Range("C12").Clear()
$ws.Range("C13").Clear()
$ws.Range("C14").Clear()
$ws.Range("C18").Clear()
$ws.Range("C19").Clear()
$ws.Range("C22").Clear()
$ws.Range("C23").Clear()
$ws.Range("C24").Clear()
$ws.Range("C25").Clear()

# J1 header connector description correction (dropped "tall")
$ws.Range("D2").Value = "20-pin 2.54mm pitch female header"

# J2 antenna cable: swap RG316 pigtail for RG174 assembly
$ws.Range("D3").Value = "RG174"
$ws.Range("F3").Value = "Antenna cable assembly"

# Narrow the now mostly-empty Status column
$ws.Columns("C").ColumnWidth = 8.8333333

# Update selection to reflect where the editor left off
$ws.Range("D26").Select()
